$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 342.19232
$ws.Range("I28").Value = 290.9565
$ws.Range("J28").Value = 735
$ws.Range("K28").Value = 290.9565
$ws.Range("L28").Value = 735
$ws.Range("M28").Value = 194.0435
$ws.Range("N28").Value = -1705

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 24867.213
$ws.Range("I98").Value = 972.0454999999999
$ws.Range("J98").Value = 72657.55
$ws.Range("K98").Value = 972.0454999999999
$ws.Range("L98").Value = 72657.55
$ws.Range("M98").Value = 525.9545000000001
$ws.Range("N98").Value = -75653.55

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 24867.213
$ws.Range("I122").Value = 972.0454999999999
$ws.Range("J122").Value = 72657.55
$ws.Range("K122").Value = 2916.1365
$ws.Range("L122").Value = 217972.65
$ws.Range("M122").Value = -466.1364999999996
$ws.Range("N122").Value = -222872.65

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1588.6428
$ws.Range("J129").Value = 1523.8096
$ws.Range("L129").Value = 4571.4288
$ws.Range("N129").Value = -14571.4288

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 36875.297
$ws.Range("I132").Value = 5265.909
$ws.Range("J132").Value = 175956.6
$ws.Range("K132").Value = 15797.727
$ws.Range("L132").Value = 527869.8
$ws.Range("M132").Value = -13267.727
$ws.Range("N132").Value = -532929.8

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3206.1035
$ws.Range("I137").Value = 1057.421
$ws.Range("J137").Value = 7288.6
$ws.Range("K137").Value = 3172.263
$ws.Range("L137").Value = 21865.8
$ws.Range("M137").Value = -622.2629999999999
$ws.Range("N137").Value = -26965.8

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1534.8596
$ws.Range("I141").Value = 1005.22644
$ws.Range("J141").Value = 8552.5
$ws.Range("K141").Value = 3015.67932
$ws.Range("L141").Value = 25657.5
$ws.Range("M141").Value = 2164.32068
$ws.Range("N141").Value = -36017.5

# ARM row 8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 1503
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 3001
$ws.Range("K8").Value = 5
$ws.Range("L8").Value = 3001
$ws.Range("M8").Value = 139
$ws.Range("N8").Value = -3289

# ARM row 10
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 35252.5
$ws.Range("J10").Value = 35252.5
$ws.Range("L10").Value = 35252.5
$ws.Range("N10").Value = -35592.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14095.081
$ws.Range("I32").Value = 11675.122
$ws.Range("K32").Value = 11675.122
$ws.Range("M32").Value = -11388.122

# ARM row 82
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 28000
$ws.Range("J82").Value = 28000
$ws.Range("L82").Value = 28000
$ws.Range("N82").Value = -28722

# ARM row 85
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H85").Value = 28000
$ws.Range("J85").Value = 28000
$ws.Range("L85").Value = 28000
$ws.Range("N85").Value = -30496

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6945868.5
$ws.Range("I132").Value = 8334300
$ws.Range("J132").Value = 3713.0833
$ws.Range("K132").Value = 25002900
$ws.Range("L132").Value = 11139.2499
$ws.Range("M132").Value = -25000370
$ws.Range("N132").Value = -16199.2499

# CRP row 12
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 20230.615
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 20230.615
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 20230.615
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -20570.615

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1007.5484
$ws.Range("I58").Value = 688.9729599999999
$ws.Range("J58").Value = 1479.04
$ws.Range("K58").Value = 688.9729599999999
$ws.Range("L58").Value = 1479.04
$ws.Range("M58").Value = -485.9729599999999
$ws.Range("N58").Value = -1885.04

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2601.375
$ws.Range("I99").Value = 2505.5
$ws.Range("J99").Value = 2633.3333
$ws.Range("K99").Value = 2505.5
$ws.Range("L99").Value = 2633.3333
$ws.Range("M99").Value = -1007.5
$ws.Range("N99").Value = -5629.3333

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2601.375
$ws.Range("I126").Value = 2505.5
$ws.Range("J126").Value = 2633.3333
$ws.Range("K126").Value = 7516.5
$ws.Range("L126").Value = 7899.999899999999
$ws.Range("M126").Value = -5046.5
$ws.Range("N126").Value = -12839.9999

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 27395.166
$ws.Range("I132").Value = 1200.1052
$ws.Range("J132").Value = 89608.44
$ws.Range("K132").Value = 3600.3156
$ws.Range("L132").Value = 268825.32
$ws.Range("M132").Value = -1070.3156
$ws.Range("N132").Value = -273885.32

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 360665.47
$ws.Range("I134").Value = 1210.9231
$ws.Range("J134").Value = 1079574.5
$ws.Range("K134").Value = 3632.7693
$ws.Range("L134").Value = 3238723.5
$ws.Range("M134").Value = -1097.7693
$ws.Range("N134").Value = -3243793.5

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1007.5484
$ws.Range("I136").Value = 688.9729599999999
$ws.Range("J136").Value = 1479.04
$ws.Range("K136").Value = 2066.91888
$ws.Range("L136").Value = 4437.12
$ws.Range("M136").Value = 483.0811200000003
$ws.Range("N136").Value = -9537.119999999999

# CUL row 94
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 7500
$ws.Range("J94").Value = 7500
$ws.Range("L94").Value = 22500
$ws.Range("N94").Value = -23852

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4855.926
$ws.Range("I131").Value = 6288.9414
$ws.Range("J131").Value = 2419.8
$ws.Range("K131").Value = 18866.8242
$ws.Range("L131").Value = 7259.400000000001
$ws.Range("M131").Value = -13826.8242
$ws.Range("N131").Value = -17339.4

# GSM row 111
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10214.714
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 11583.833
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 11583.833
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -11959.833

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1948.909
$ws.Range("I132").Value = 820.1
$ws.Range("J132").Value = 2889.5833
$ws.Range("K132").Value = 2460.3
$ws.Range("L132").Value = 8668.749899999999
$ws.Range("M132").Value = 69.69999999999982
$ws.Range("N132").Value = -13728.7499

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 18410.807
$ws.Range("I136").Value = 22765.467
$ws.Range("J136").Value = 2080.8333
$ws.Range("K136").Value = 68296.401
$ws.Range("L136").Value = 6242.499899999999
$ws.Range("M136").Value = -65746.401
$ws.Range("N136").Value = -11342.4999
